# Weekly data refresh: a new week's price record (2021-09-22) is added for
# "Macroferia Regional de Talca - Repollo" in front of the existing history,
# pushing the prior rows (159-170) down by one (to 160-171).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 159, shifting rows 159:170 down to 160:171.
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with this week's record.
$ws.Cells.Item(159, 1).Value = 5
$ws.Cells.Item(159, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(159, 3).Value = "Maule"
$ws.Cells.Item(159, 4).Value = "2021-09-22"
$ws.Cells.Item(159, 5).Value = 7
$ws.Cells.Item(159, 6).Value = 100112006
$ws.Cells.Item(159, 7).Value = "Repollo"
$ws.Cells.Item(159, 8).Value = "Crespo record"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 3000
$ws.Cells.Item(159, 11).Value = 500
$ws.Cells.Item(159, 12).Value = 500
$ws.Cells.Item(159, 13).Value = 500
$ws.Cells.Item(159, 14).Value = "$/unidad"
$ws.Cells.Item(159, 15).Value = "Región del Maule"
$ws.Cells.Item(159, 16).Value = 500
$ws.Cells.Item(159, 17).Value = 1
$ws.Cells.Item(159, 18).Value = "Hortaliza"
